$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2096.9473
$ws.Range("J62").Value = 2027.1111
$ws.Range("L62").Value = 2027.1111
$ws.Range("N62").Value = -3275.1111
# Row 65
$ws.Range("H65").Value = 2096.9473
$ws.Range("J65").Value = 2027.1111
$ws.Range("L65").Value = 10135.5555
$ws.Range("N65").Value = -16375.5555
# Row 106
$ws.Range("H106").Value = 1367.6364
$ws.Range("I106").Value = 1100.1111
$ws.Range("J106").Value = 2571.5
$ws.Range("K106").Value = 1100.1111
$ws.Range("L106").Value = 2571.5
$ws.Range("M106").Value = -469.1111000000001
$ws.Range("N106").Value = -3833.5
# Row 113
$ws.Range("H113").Value = 20003900
$ws.Range("I113").Value = 100002150
$ws.Range("J113").Value = 4337.175
$ws.Range("K113").Value = 100002150
$ws.Range("L113").Value = 4337.175
$ws.Range("M113").Value = -99998896
$ws.Range("N113").Value = -10845.175
# Row 115
$ws.Range("H115").Value = 674
$ws.Range("J115").Value = 600
$ws.Range("L115").Value = 1800
$ws.Range("N115").Value = -4934
# Row 129
$ws.Range("H129").Value = 616
# Row 132
$ws.Range("H132").Value = 41651.08
$ws.Range("I132").Value = 51213.57
$ws.Range("J132").Value = 1488.6
$ws.Range("K132").Value = 153640.71
$ws.Range("L132").Value = 4465.799999999999
$ws.Range("M132").Value = -151110.71
$ws.Range("N132").Value = -9525.799999999999
# Row 135
$ws.Range("H135").Value = 6944.9585
$ws.Range("I135").Value = 1474.5
$ws.Range("J135").Value = 10852.429
$ws.Range("K135").Value = 13270.5
$ws.Range("L135").Value = 97671.861
$ws.Range("M135").Value = -10735.5
$ws.Range("N135").Value = -102741.861
# Row 137
$ws.Range("H137").Value = 29558.73
$ws.Range("I137").Value = 2818.5
$ws.Range("J137").Value = 61017.824
$ws.Range("K137").Value = 8455.5
$ws.Range("L137").Value = 183053.472
$ws.Range("M137").Value = -5905.5
$ws.Range("N137").Value = -188153.472
# Row 138
$ws.Range("H138").Value = 3480.5
$ws.Range("I138").Value = 1665.6666
$ws.Range("J138").Value = 3739.762
$ws.Range("K138").Value = 4996.9998
$ws.Range("L138").Value = 11219.286
$ws.Range("M138").Value = 143.0002000000004
$ws.Range("N138").Value = -21499.286

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1198.7755
$ws.Range("I2").Value = 959.2778
$ws.Range("J2").Value = 1862
$ws.Range("K2").Value = 959.2778
$ws.Range("L2").Value = 1862
$ws.Range("M2").Value = -846.2778
$ws.Range("N2").Value = -2088
# Row 32
$ws.Range("H32").Value = 19555.732
$ws.Range("I32").Value = 24436.133
$ws.Range("K32").Value = 24436.133
$ws.Range("M32").Value = -24149.133
# Row 61
$ws.Range("H61").Value = 312930.44
$ws.Range("I61").Value = 451993.78
$ws.Range("J61").Value = 3900.7778
$ws.Range("K61").Value = 451993.78
$ws.Range("L61").Value = 3900.7778
$ws.Range("M61").Value = -451781.78
$ws.Range("N61").Value = -4324.7778
# Row 74
$ws.Range("H74").Value = 2221.081
$ws.Range("I74").Value = 2737.05
$ws.Range("J74").Value = 1614.0588
$ws.Range("K74").Value = 2737.05
$ws.Range("L74").Value = 1614.0588
$ws.Range("M74").Value = -1863.05
$ws.Range("N74").Value = -3362.0588
# Row 77
$ws.Range("H77").Value = 2221.081
$ws.Range("I77").Value = 2737.05
$ws.Range("J77").Value = 1614.0588
$ws.Range("K77").Value = 13685.25
$ws.Range("L77").Value = 8070.294
$ws.Range("M77").Value = -9317.25
$ws.Range("N77").Value = -16806.294
# Row 97
$ws.Range("H97").Value = 2367.889
$ws.Range("J97").Value = 1798.2
$ws.Range("L97").Value = 1798.2
$ws.Range("N97").Value = -2790.2
# Row 116
$ws.Range("H116").Value = 1198.7755
$ws.Range("I116").Value = 959.2778
$ws.Range("J116").Value = 1862
$ws.Range("K116").Value = 959.2778
$ws.Range("L116").Value = 1862
$ws.Range("M116").Value = 1334.7222
$ws.Range("N116").Value = -6450
# Row 136
$ws.Range("H136").Value = 312930.44
$ws.Range("I136").Value = 451993.78
$ws.Range("J136").Value = 3900.7778
$ws.Range("K136").Value = 1355981.34
$ws.Range("L136").Value = 11702.3334
$ws.Range("M136").Value = -1353431.34
$ws.Range("N136").Value = -16802.3334
# Row 139
$ws.Range("H139").Value = 50857.31
$ws.Range("J139").Value = 50857.31
$ws.Range("L139").Value = 50857.31
$ws.Range("N139").Value = -61137.31

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1198.7755
$ws.Range("I3").Value = 959.2778
$ws.Range("J3").Value = 1862
$ws.Range("K3").Value = 959.2778
$ws.Range("L3").Value = 1862
$ws.Range("M3").Value = -845.2778
$ws.Range("N3").Value = -2090
# Row 134
$ws.Range("H134").Value = 32787.695
$ws.Range("I134").Value = 38475.32
$ws.Range("J134").Value = 937
$ws.Range("K134").Value = 115425.96
$ws.Range("L134").Value = 2811
$ws.Range("M134").Value = -112890.96
$ws.Range("N134").Value = -7881

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10675.94
$ws.Range("I31").Value = 17968.16
$ws.Range("J31").Value = 3383.72
$ws.Range("K31").Value = 17968.16
$ws.Range("L31").Value = 3383.72
$ws.Range("M31").Value = -17673.16
$ws.Range("N31").Value = -3973.72
# Row 34
$ws.Range("H34").Value = 10675.94
$ws.Range("I34").Value = 17968.16
$ws.Range("J34").Value = 3383.72
$ws.Range("K34").Value = 17968.16
$ws.Range("L34").Value = 3383.72
$ws.Range("M34").Value = -17766.16
$ws.Range("N34").Value = -3787.72

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 2949.75
$ws.Range("I80").Value = 2949.5
$ws.Range("J80").Value = 2950
$ws.Range("K80").Value = 8848.5
$ws.Range("L80").Value = 8850
$ws.Range("M80").Value = -7912.5
$ws.Range("N80").Value = -10722
# Row 83
$ws.Range("H83").Value = 2949.75
$ws.Range("I83").Value = 2949.5
$ws.Range("J83").Value = 2950
$ws.Range("K83").Value = 26545.5
$ws.Range("L83").Value = 26550
$ws.Range("M83").Value = -21865.5
$ws.Range("N83").Value = -35910
# Row 107
$ws.Range("H107").Value = 4673
$ws.Range("J107").Value = 1085.6957
$ws.Range("L107").Value = 3257.0871
$ws.Range("N107").Value = -7097.0871
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
# Row 131
$ws.Range("H131").Value = 143685.88
$ws.Range("J131").Value = 162129.38
$ws.Range("L131").Value = 486388.14
$ws.Range("N131").Value = -496468.14

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 42
$ws.Range("H42").Value = 31290
$ws.Range("J42").Value = 31290
$ws.Range("L42").Value = 31290
$ws.Range("N42").Value = -32260
# Row 115
$ws.Range("H115").Value = 31290
$ws.Range("J115").Value = 31290
$ws.Range("L115").Value = 31290
$ws.Range("N115").Value = -33640
# Row 132
$ws.Range("H132").Value = 39991.824
$ws.Range("I132").Value = 55601.895
$ws.Range("J132").Value = 25868.428
$ws.Range("K132").Value = 166805.685
$ws.Range("L132").Value = 77605.284
$ws.Range("M132").Value = -164275.685
$ws.Range("N132").Value = -82665.284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3427.75
$ws.Range("I68").Value = 1905.7778
$ws.Range("J68").Value = 4673
$ws.Range("K68").Value = 1905.7778
$ws.Range("L68").Value = 4673
$ws.Range("M68").Value = -1156.7778
$ws.Range("N68").Value = -6171
# Row 71
$ws.Range("H71").Value = 3427.75
$ws.Range("I71").Value = 1905.7778
$ws.Range("J71").Value = 4673
$ws.Range("K71").Value = 9528.889000000001
$ws.Range("L71").Value = 23365
$ws.Range("M71").Value = -5784.889000000001
$ws.Range("N71").Value = -30853
# Row 132
$ws.Range("H132").Value = 1448.6078
$ws.Range("I132").Value = 1025.9736
$ws.Range("K132").Value = 3077.9208
$ws.Range("M132").Value = -547.9207999999999
# Row 141
$ws.Range("H141").Value = 56333.332
$ws.Range("J141").Value = 56333.332
$ws.Range("L141").Value = 56333.332
$ws.Range("N141").Value = -66693.33199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1610.7727
$ws.Range("I136").Value = 983.3
$ws.Range("J136").Value = 2133.6667
$ws.Range("K136").Value = 2949.9
$ws.Range("L136").Value = 6401.000100000001
$ws.Range("M136").Value = -399.8999999999996
$ws.Range("N136").Value = -11501.0001
